# Update the handback report timestamps, as if re-generated a little later.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 9434bfbc-...md
$wsOverview.Range("G2").Value = "2016-08-26 09:11:46"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 9434bfbc-...md
$wsZhCn.Range("H2").Value = "2016-08-26 09:11:42"
$wsZhCn.Range("K2").Value = "2016-08-26 09:11:59"

# de-de sheet: "Latest HO Xliff Generate Date" (shared with Overview) and "Correspond Handback DateTime"
$wsDeDe.Range("H2").Value = "2016-08-26 09:11:46"
$wsDeDe.Range("K2").Value = "2016-08-26 09:12:16"
